# "Wrote script for product search"
# Adds a new Sheet2 (after Sheet1) holding a small product-search table,
# and makes it the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet immediately after Sheet1 (Add() defaults to "before
# the active sheet", so pass Sheet1 explicitly as the After target).
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header/product rows.
$ws2.Range("A1").Value = "OnePlus 8 (Glacial Green 6GB RAM+128GB Storage)"
$ws2.Range("B1").Value = "Electronics"
$ws2.Range("C1").Value = 39999
$ws2.Range("D1").Value = 4.3

$ws2.Range("A2").Value = "OnePlus 8 Pro (Onyx Black 8GB RAM+128GB Storage)"
$ws2.Range("B2").Value = "Electronics"
$ws2.Range("C2").Value = 54999
$ws2.Range("D2").Value = 4.1

# Widen column A so the product names are readable.
$ws2.Columns.Item(1).ColumnWidth = 61.5

# Sheet2 becomes the active sheet/tab, with its own selection at E10.
$ws2.Range("E10").Select()
$ws2.Activate()
